$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting (values such as "1.00" or
# "39.870.72" must not be auto-converted to numbers by Excel when assigned via .Value)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.870.72"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.206.23"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "288.33"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "87.12"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "30.33"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "0.0776"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").Value = "6.42"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "2.550.41"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "13.88"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "2.207.61"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "39.803.65"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "11.54"
$ws.Range("E19").Value = "  +9.68%  "
$ws.Range("D20").Value = "0.0₃0879"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "5.77"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "65.32"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "235.23"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "1.81"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "22.42"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "152.62"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").Value = "31.62"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "4.92"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").Value = "0.0715"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "2.80"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "15.66"
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.0984"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "1.69"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +3.73%  "
$ws.Range("D42").Value = "2.088.64"
$ws.Range("E42").Value = "  +7.30%  "
$ws.Range("D43").Value = "2.14"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "9.94"
$ws.Range("E44").Value = "  +5.78%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0267"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "17.50"
$ws.Range("E46").Value = "  +7.35%  "
$ws.Range("D47").Value = "2.64"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "2.423.50"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "1.44"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "88.20"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "68.74"
$ws.Range("E51").Value = "  -2.90%  "
